$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns for the rows we touch are treated as plain text,
# so numeric-looking strings (e.g. "1.001", "0.9997") are preserved verbatim rather than
# being auto-converted into numbers by Excel.
$textRanges = @("D2:D19", "D21:D23", "D25:D33", "D35:D46", "D48:D51", "E2:E9", "E12:E17", "E19:E51")
foreach ($addr in $textRanges) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.483.81"
$ws.Range("E2").Value = "  -5.24%  "

$ws.Range("D3").Value = "1.838.85"
$ws.Range("E3").Value = "  -4.26%  "

$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.52%  "

$ws.Range("D5").Value = "313.08"
$ws.Range("E5").Value = "  -3.72%  "

$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("D7").Value = "0.4229"
$ws.Range("E7").Value = "  -7.82%  "

$ws.Range("D8").Value = "0.3627"
$ws.Range("E8").Value = "  -4.78%  "

$ws.Range("D9").Value = "44.30"
$ws.Range("E9").Value = "  -2.95%  "

$ws.Range("D10").Value = "0.07219"

$ws.Range("D11").Value = "0.9036"

$ws.Range("D12").Value = "20.59"
$ws.Range("E12").Value = "  -8.98%  "

$ws.Range("D13").Value = "1.813.41"
$ws.Range("E13").Value = "  -5.99%  "

$ws.Range("D14").Value = "6.566"
$ws.Range("E14").Value = "  -5.43%  "

$ws.Range("D15").Value = "5.331"
$ws.Range("E15").Value = "  -6.63%  "

$ws.Range("D16").Value = "0.06798"
$ws.Range("E16").Value = "  -3.11%  "

$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").Value = "77.45"

$ws.Range("D19").Value = "0.000008942"
$ws.Range("E19").Value = "  -5.81%  "

$ws.Range("E20").Value = "  -0.36%  "

$ws.Range("D21").Value = "15.35"
$ws.Range("E21").Value = "  -7.95%  "

$ws.Range("D22").Value = "27.497.83"
$ws.Range("E22").Value = "  -5.26%  "

$ws.Range("D23").Value = "4.934"
$ws.Range("E23").Value = "  -7.75%  "

$ws.Range("E24").Value = "  -4.46%  "

$ws.Range("D25").Value = "1.995.01"
$ws.Range("E25").Value = "  -7.41%  "

$ws.Range("D26").Value = "2.013"
$ws.Range("E26").Value = "  -2.34%  "

$ws.Range("D27").Value = "152.96"
$ws.Range("E27").Value = "  -3.11%  "

$ws.Range("D28").Value = "18.15"
$ws.Range("E28").Value = "  -4.68%  "

$ws.Range("D29").Value = "5.250"
$ws.Range("E29").Value = "  -6.07%  "

$ws.Range("D30").Value = "110.65"
$ws.Range("E30").Value = "  -5.73%  "

$ws.Range("D31").Value = "1.665"
$ws.Range("E31").Value = "  -8.97%  "

$ws.Range("D32").Value = "0.08862"
$ws.Range("E32").Value = "  -4.95%  "

$ws.Range("D33").Value = "0.7767"
$ws.Range("E33").Value = "  -9.47%  "

$ws.Range("E34").Value = "  -11.28%  "

$ws.Range("D35").Value = "2.896"
$ws.Range("E35").Value = "  -4.17%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.074"
$ws.Range("E36").Value = "  -13.24%  "

$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "0.9996"
$ws.Range("E37").Value = "  -0.46%  "

$ws.Range("D38").Value = "0.05334"
$ws.Range("E38").Value = "  -6.08%  "

$ws.Range("D39").Value = "1.076"
$ws.Range("E39").Value = "  -6.49%  "

$ws.Range("D40").Value = "0.01927"
$ws.Range("E40").Value = "  -5.51%  "

$ws.Range("D41").Value = "2.936"
$ws.Range("E41").Value = "  -5.80%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.5055"
$ws.Range("E42").Value = "  -7.86%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "6.810"
$ws.Range("E43").Value = "  -8.23%  "

$ws.Range("D44").Value = "0.1632"
$ws.Range("E44").Value = "  -6.94%  "

$ws.Range("D45").Value = "0.06618"
$ws.Range("E45").Value = "  -4.32%  "

$ws.Range("D46").Value = "8.220"
$ws.Range("E46").Value = "  -12.01%  "

$ws.Range("E47").Value = "  -8.91%  "

$ws.Range("D48").Value = "105.08"
$ws.Range("E48").Value = "  -4.64%  "

$ws.Range("D49").Value = "10.18"
$ws.Range("E49").Value = "  -8.91%  "

$ws.Range("D50").Value = "0.9986"
$ws.Range("E50").Value = "  -0.55%  "

$ws.Range("D51").Value = "1.628"
$ws.Range("E51").Value = "  -7.27%  "
